# Auto-generated edit script: adds 20 new quiz rows (101-120) to Sheet1
# with their shared strings, plus the associated cell styles (wrap text
# on multi-line question/answer cells, autofit-equivalent row heights,
# and the distinct font style observed on cell C107).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 101
$ws.Range("A101").Value = "Pytanie:  W jednostkach PSP jako podstawowy system zabezpieczenia dróg
oddechowych stosuje się:"
$ws.Range("B101").Value = "aparaty tlenowe"
$ws.Range("C101").Value = "aparaty powietrzne nadciśnieniowe"
$ws.Range("D101").Value = "aparaty powietrzne podciśnieniowe"
$ws.Range("F101").Value = 2
$ws.Range("A101").WrapText = $true
$ws.Rows.Item(101).RowHeight = 30

# Row 102
$ws.Range("A102").Value = "Pytanie: W jakim stopniu jest obecny Komendant Wojewódzki Państwowej
Straży Pożarnej w Warszawie?"
$ws.Range("B102").Value = "starszy brygadie"
$ws.Range("C102").Value = "nadbrygadier generał"
$ws.Range("D102").Value = "nadbrygadier"
$ws.Range("F102").Value = 3
$ws.Range("A102").WrapText = $true
$ws.Rows.Item(102).RowHeight = 30

# Row 103
$ws.Range("A103").Value = "Pytanie: Flame Control to środek"
$ws.Range("B103").Value = "gaśniczy"
$ws.Range("C103").Value = "ogniochronny"
$ws.Range("D103").Value = "absorbent"
$ws.Range("F103").Value = 2

# Row 104
$ws.Range("A104").Value = "Pytanie: Prawo do wyznaczania w lesie miejsca na palenie ogniska ma:"
$ws.Range("B104").Value = "właściciel lub zarządca lasu"
$ws.Range("C104").Value = "Komendant Powiatowy lub Miejski PSP"
$ws.Range("D104").Value = "osoba dorosła, pod warunkiem zachowania środków ostrożności"
$ws.Range("F104").Value = 1

# Row 105
$ws.Range("A105").Value = "Pytanie: Do gazów gaśniczych można zaliczyć:"
$ws.Range("B105").Value = "azot"
$ws.Range("C105").Value = "etan"
$ws.Range("D105").Value = "metan"
$ws.Range("F105").Value = 1

# Row 106
$ws.Range("A106").Value = "Pytanie: Pianą można gasić pożary klasy:"
$ws.Range("B106").Value = "klasy A i B"
$ws.Range("C106").Value = "klasy C i D"
$ws.Range("D106").Value = "tylko klasy C"
$ws.Range("F106").Value = 1

# Row 107
$ws.Range("A107").Value = "Pytanie:  Pożary metali lekkich zaliczają się do klasy pożarów:"
$ws.Range("B107").Value = "A"
$ws.Range("C107").Value = "B"
$ws.Range("D107").Value = "D"
$ws.Range("F107").Value = 3
$ws.Range("C107").Font.Name = "Calibri"

# Row 108
$ws.Range("A108").Value = "Pytanie: Po haśle „Do akcji gotuj” przy zasilaniu samochodu z hydrantu,
kierowca:"
$ws.Range("B108").Value = "pomaga w budowie linii głównej"
$ws.Range("C108").Value = "przeprowadza rozpoznanie wodne"
$ws.Range("D108").Value = "ustawia samochód w miejscu wskazanym przez dowódcę, a następnie uruchamia autopompę"
$ws.Range("F108").Value = 3
$ws.Range("A108").WrapText = $true
$ws.Rows.Item(108).RowHeight = 30

# Row 109
$ws.Range("A109").Value = "Pytanie: Czego nie wolno gasić wodą?"
$ws.Range("B109").Value = "drewna"
$ws.Range("C109").Value = "węgla"
$ws.Range("D109").Value = "sodu"
$ws.Range("F109").Value = 3

# Row 110
$ws.Range("A110").Value = "Pytanie: Jaki narząd jest najbardziej wrażliwy na niedotlenienie?"
$ws.Range("B110").Value = "serce"
$ws.Range("C110").Value = "mózg"
$ws.Range("D110").Value = "nerki"
$ws.Range("F110").Value = 2

# Row 111
$ws.Range("A111").Value = "Pytanie: W mięśniu czwórgłowym uda znajduje się wbity pręt metalowy
pierwsza pomoc polega na:"
$ws.Range("B111").Value = "szybkim usunięciu przedmiotu"
$ws.Range("C111").Value = "zabezpieczenie przedmiotu przed przemieszczeniem i opatrzenie
rany"
$ws.Range("D111").Value = "zebraniu danych osobowych póki pacjent jest przytomny"
$ws.Range("F111").Value = 2
$ws.Range("A111").WrapText = $true
$ws.Range("C111").WrapText = $true
$ws.Rows.Item(111).RowHeight = 30

# Row 112
$ws.Range("A112").Value = "Pytanie: Wstrząs jest to:"
$ws.Range("B112").Value = "stan niedotlenienia mózgu"
$ws.Range("C112").Value = "zespół drgawek"
$ws.Range("D112").Value = " odruch człowieka na działanie prądu elektrycznego"
$ws.Range("F112").Value = 1

# Row 113
$ws.Range("A113").Value = "Pytanie: Krew tętniczą wypływająca z rany poznamy po:"
$ws.Range("B113").Value = "ciemnobrunatnej barwie"
$ws.Range("C113").Value = "spokojnym wypływie"
$ws.Range("D113").Value = "jasnoczerwonej barwie"
$ws.Range("F113").Value = 3

# Row 114
$ws.Range("A114").Value = "Pytanie: Co oznacza znak gestowy „Prawa ręka podniesiona w górę”:"
$ws.Range("B114").Value = "woda stój"
$ws.Range("C114").Value = "woda naprzód"
$ws.Range("D114").Value = "uwaga"
$ws.Range("F114").Value = 3

# Row 115
$ws.Range("A115").Value = "Pytanie:  Mostek przejazdowy służy do:"
$ws.Range("B115").Value = "umożliwienia przejazdu samochodu przez rzekę"
$ws.Range("C115").Value = "zabezpieczenia węży przed uszkodzeniem układanych na ulicach i
drogach"
$ws.Range("D115").Value = "sprawiania drabiny pożarniczej"
$ws.Range("F115").Value = 2
$ws.Range("C115").WrapText = $true
$ws.Rows.Item(115).RowHeight = 30

# Row 116
$ws.Range("A116").Value = "Pytanie: Dym biało-żółty wydziela się podczas palenia:"
$ws.Range("B116").Value = "siana"
$ws.Range("C116").Value = "celuloidu"
$ws.Range("D116").Value = "gumy"
$ws.Range("F116").Value = 1

# Row 117
$ws.Range("A117").Value = "Pytanie: Stosunek służbowy strażaka powstaje z dniem:"
$ws.Range("B117").Value = "powołania"
$ws.Range("C117").Value = "mianowania"
$ws.Range("D117").Value = "ślubowania"
$ws.Range("F117").Value = 2

# Row 118
$ws.Range("A118").Value = "Pytanie: Podoficerom pełniącym służbę w jednostkach ratowniczo –
gaśniczych kolejne stopnie służbowe nadaje:"
$ws.Range("B118").Value = " Komendant /Miejski/ Powiatowy PSP"
$ws.Range("C118").Value = "Komendant Wojewódzki PSP"
$ws.Range("D118").Value = "Komendant Główny PSP"
$ws.Range("F118").Value = 2
$ws.Range("A118").WrapText = $true
$ws.Rows.Item(118).RowHeight = 30

# Row 119
$ws.Range("A119").Value = "Pytanie: Wytwornica pianowa WP 2/75 służy do wytwarzania piany:"
$ws.Range("B119").Value = "średniej"
$ws.Range("C119").Value = "lekkiej"
$ws.Range("D119").Value = "ciężkiej"
$ws.Range("F119").Value = 1

# Row 120
$ws.Range("A120").Value = "Pytanie: W jakim mieści na terenie woj. kujawsko pomorskiego znajduje się
Szkoła Podoficerska?"
$ws.Range("B120").Value = "Toruń"
$ws.Range("C120").Value = "Włocławek"
$ws.Range("D120").Value = "Bydgoszcz"
$ws.Range("F120").Value = 3
$ws.Range("A120").WrapText = $true
$ws.Rows.Item(120).RowHeight = 30

# Restore the last-used selection noted in the source workbook.
$ws.Range("E124").Select()

